# Edit script generated to apply the weekly crime-data refresh described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Type-changing cells: numeric <-> dash ("-") / "***.*" placeholders ---
# Row 23 (Housing) has untouched template cells we copy style+value from:
#   C23/D23/F23/G23 = dash ("-", style 14, shared string index 20)
#   E23/H23         = "***.*" (style 14, shared string index 21)

function Set-NumberCell($addr, $value, $format) {
    $ws.Range($addr).Value = $value
    $ws.Range($addr).NumberFormat = $format
}

function Set-DashCell($addr, $templateAddr) {
    $ws.Range($templateAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($templateAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

# Row 15 (Rape)
Set-NumberCell "D15" 1 "#,##0"
Set-NumberCell "E15" -100 "#,##0.0;""-""#,##0.0"
Set-NumberCell "G15" 1 "#,##0"
Set-NumberCell "H15" -100 "#,##0.0;""-""#,##0.0"

# Row 22 (Transit)
Set-DashCell "C22" "C23"

# Row 26 (UCR Rape*)
Set-NumberCell "D26" 1 "#,##0"
Set-NumberCell "E26" -100 "#,##0.0;""-""#,##0.0"
Set-NumberCell "G26" 1 "#,##0"
Set-NumberCell "H26" -100 "#,##0.0;""-""#,##0.0"

# Row 27 (Other Sex Crimes)
Set-DashCell "C27" "C23"
Set-NumberCell "D27" 1 "#,##0"
Set-NumberCell "E27" -100 "#,##0.0;""-""#,##0.0"

# --- Plain numeric value refreshes (style/format unchanged) ---
# Row 15
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 46.153846153846
$ws.Range("N15").Value = -40.625

# Row 16
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -71.428571428571
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -44
$ws.Range("I16").Value = 220
$ws.Range("J16").Value = 199
$ws.Range("K16").Value = 10.552763819095
$ws.Range("L16").Value = 25
$ws.Range("M16").Value = -16.349809885931
$ws.Range("N16").Value = -74.828375286041

# Row 17
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 25
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -10
$ws.Range("I17").Value = 337
$ws.Range("J17").Value = 338
$ws.Range("K17").Value = -0.295857988165
$ws.Range("L17").Value = 21.660649819494
$ws.Range("M17").Value = 126.174496644295
$ws.Range("N17").Value = -16.37717121588

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -62.5
$ws.Range("I18").Value = 143
$ws.Range("J18").Value = 124
$ws.Range("K18").Value = 15.322580645161
$ws.Range("L18").Value = -11.180124223602
$ws.Range("M18").Value = -52.960526315789
$ws.Range("N18").Value = -88.288288288288

# Row 19
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 75
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = 31.111111111111
$ws.Range("I19").Value = 656
$ws.Range("J19").Value = 509
$ws.Range("K19").Value = 28.880157170923
$ws.Range("L19").Value = 26.396917148362
$ws.Range("M19").Value = 79.72602739726
$ws.Range("N19").Value = 11.186440677966

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 218
$ws.Range("J20").Value = 211
$ws.Range("K20").Value = 3.317535545023
$ws.Range("L20").Value = 19.125683060109
$ws.Range("M20").Value = -32.50773993808
$ws.Range("N20").Value = -93.302611367127

# Row 21
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -13.513513513513
$ws.Range("F21").Value = 116
$ws.Range("G21").Value = 136
$ws.Range("H21").Value = -14.705882352941
$ws.Range("I21").Value = 1598
$ws.Range("J21").Value = 1399
$ws.Range("K21").Value = 14.22444603288
$ws.Range("L21").Value = 20.240782543265
$ws.Range("M21").Value = 12.061711079943
$ws.Range("N21").Value = -75.015634771732

# Row 24
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -12
$ws.Range("F24").Value = 113
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = 15.306122448979
$ws.Range("I24").Value = 1345
$ws.Range("J24").Value = 978
$ws.Range("K24").Value = 37.525562372188
$ws.Range("L24").Value = 60.119047619047
$ws.Range("M24").Value = 97.794117647058

# Row 25
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = -21.951219512195
$ws.Range("I25").Value = 520
$ws.Range("J25").Value = 466
$ws.Range("K25").Value = 11.587982832618
$ws.Range("L25").Value = 18.993135011441
$ws.Range("M25").Value = 7.883817427385

# Row 26
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = 84.210526315789

# Row 27
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -33.333333333333
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 40
$ws.Range("L27").Value = 80.645161290322

$excel.CutCopyMode = $false

